# Add a new "2022-Q3" quarter: insert a new sheet with its fund-holding
# detail data, and add a corresponding summary row at the top of "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q3,
#    shifting the existing quarters down by one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Shift rows 2..7 down to 3..8 (iterate bottom-up so we don't clobber data).
for ($r = 7; $r -ge 2; $r--) {
    $newR = $r + 1
    $summary.Cells.Item($newR, 2).Value = $summary.Cells.Item($r, 2).Value()
    $summary.Cells.Item($newR, 3).Value = $summary.Cells.Item($r, 3).Value()
    $summary.Cells.Item($newR, 4).Value = $summary.Cells.Item($r, 4).Value()
}

# Write the new 2022-Q3 summary row.
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 17
$summary.Cells.Item(2, 4).Value = 4.06

# Column A is just a 0-based row counter; rewrite it for every data row.
for ($r = 2; $r -le 8; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# Row 8 is brand new (the sheet only had 7 rows before), so its column-A
# cell doesn't yet carry the bold/centered/bordered style used by the rest
# of that column; copy it over from a cell that already has it.
$summary.Range("A2").Copy()
$summary.Range("A8").PasteSpecial(-4122)
$summary.Range("A8").Value = 6

# ---------------------------------------------------------------------
# 2. Insert a brand-new worksheet "2022-Q3" right after "总计" (i.e.
#    before the old first quarter sheet) holding the per-fund detail.
# ---------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $firstSheet)
$q3.Name = "2022-Q3"

$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

$rows = @(
    @(0,  "310358", "申万菱信新经济混合",                 "26.05", "90.23", "7.28", "1.8964", 1),
    @(1,  "011488", "申万菱信乐享混合",                   "9.63",  "89.62", "7.02", "0.6760", 2),
    @(2,  "013085", "申万菱信乐同混合型证券投资基金A",     "9.68",  "87.21", "5.67", "0.5489", 4),
    @(3,  "012051", "申万菱信乐道三年持有期混合",           "3.39",  "90.38", "6.50", "0.2204", 2),
    @(4,  "012210", "申万菱信智能汽车股票A",               "2.98",  "93.85", "7.37", "0.2196", 2),
    @(5,  "013634", "申万菱信双利混合A",                   "6.76",  "26.09", "1.71", "0.1156", 7),
    @(6,  "200010", "长城双动力混合A",                     "3.29",  "93.10", "2.80", "0.0921", 10),
    @(7,  "012211", "申万菱信智能汽车股票C",               "1.04",  "93.85", "7.37", "0.0766", 2),
    @(8,  "015561", "长城双动力混合C",                     "2.72",  "93.10", "2.80", "0.0762", 10),
    @(9,  "013086", "申万菱信乐同混合型证券投资基金C",     "1.10",  "87.21", "5.67", "0.0624", 4),
    @(10, "001707", "诺安高端制造股票A",                   "1.12",  "90.57", "3.99", "0.0447", 6),
    @(11, "010857", "宝盈祥乐一年持有期混合型证券投资基金A", "1.08",  "27.77", "1.30", "0.0140", 10),
    @(12, "013635", "申万菱信双利混合C",                   "0.52",  "26.09", "1.71", "0.0089", 7),
    @(13, "014246", "大摩现代服务业混合A",                 "0.17",  "66.96", "2.68", "0.0046", 10),
    @(14, "014247", "大摩现代服务业混合C",                 "0.06",  "66.96", "2.68", "0.0016", 10),
    @(15, "010858", "宝盈祥乐一年持有期混合型证券投资基金C", "0.06",  "27.77", "1.30", "0.0008", 10),
    @(16, "014536", "诺安高端制造股票C",                   "0.00",  "90.57", "3.99", "0",      6)
)

$lastDataRow = 1 + $rows.Length

# Columns B (fund code), D, E, F and (for all but the last row) G must be
# stored as TEXT, not numbers, so leading zeros / exact formatting survive
# (matches the source data, which keeps these as text everywhere else too).
$q3.Range("B2:B" + $lastDataRow).NumberFormat = "@"
$q3.Range("D2:F" + $lastDataRow).NumberFormat = "@"
$q3.Range("G2:G" + ($lastDataRow - 1)).NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
}

# The very last row's "持有市值" (G) is the literal number 0, unlike the
# other rows in that column which are text.
$q3.Range("G" + $lastDataRow).NumberFormat = "General"
$q3.Range("G" + $lastDataRow).Value = 0

# Drop the forced number formats again now that the text is committed, so
# the cells end up with no explicit style (matching the rest of the file).
$q3.Range("B2:B" + $lastDataRow).Style = "Normal"
$q3.Range("D2:F" + $lastDataRow).Style = "Normal"
$q3.Range("G2:G" + $lastDataRow).Style = "Normal"

# Column A and the header row use the bold/centered/bordered style that's
# already used for the same roles elsewhere in the workbook; copy it over.
$summary.Range("A2").Copy()
$q3.Range("A2:A" + $lastDataRow).PasteSpecial(-4122)

$summary.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

Write-Host "done"
